$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) info sheet: update the single numeric result value
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("info")
$wsInfo.Range("B2").Value = 8.786959409713745

# ---------------------------------------------------------------------------
# 2) Swap the "x" and "u" sheet tab names (use a temporary name to avoid a
#    naming collision while the swap is in progress). sheetId stays bound to
#    the underlying sheet (rId), so renaming is all that is needed to match
#    <sheet name="u" sheetId="3" .../><sheet name="x" sheetId="4" .../>.
# ---------------------------------------------------------------------------
$wsOldX = $wb.Worksheets.Item("x")
$wsOldU = $wb.Worksheets.Item("u")

$wsOldX.Name = "__tmp_swap__"
$wsOldU.Name = "x"
$wsOldX.Name = "u"

# After the rename:
#   $wsOldX -> tab now named "u" (this used to be the "x" data sheet, sheetId 3)
#   $wsOldU -> tab now named "x" (this used to be the "u" data sheet, sheetId 4)
$wsU = $wsOldX
$wsX = $wsOldU

# ---------------------------------------------------------------------------
# 3) "u" sheet (tab now named "u"): add the "i" header in A1 (copying the
#    bold/boxed header style used throughout this workbook), relabel the
#    value header in B1 as "u" and write the new u-values.
# ---------------------------------------------------------------------------
$wsU.Range("A1").Value = "i"
$wsU.Range("A2").Copy()
$wsU.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$wsU.Range("B1").Value = "u"

$uVals = @(-1.025491719235965,-1.758836587582003,-0.8176983786983505,-1.693327070700814,-1.427232872661531,-1.85216118651203,-1.587089830666576,-1.131123663615027,-1.576478257915759,-0.1330505013343688,0.3619614908570461,-0.6733534270026644,0.04279391483330208,-0.2860727774270297,0.4211272049439394,-0.5516255859351369,-0.1786750279326255,-0.1711120833397799,0.1958218876858941,0.3259246801277786,0.05457021758864933,-1.997473315719062,-1.444917422460285,-1.426047813317746,-1.337645520118827,0.1921953822923932,-0.9197044060225648,-1.946130505313933,0.3171648728620471,-1.02822103442286,-1.35534069673737,-1.285006599975327,-0.8866338056867975,-1.985852824287604,-1.062808922768782,0.4810717520620247,-0.952175413440536,-1.457521079667401,0.2953801924299815,-0.7404319922607412,0.06954245483389832,-1.250452056972108,-1.864328467139036,0.3581224081269534,-1.357083759538127,-1.389442019938829,0.8261234534512258,0.3944403781887384,-0.03015841203012881,-0.9947082869631791,-0.4302660392933275,0.2936718543552272,-0.2194355189449402,-1.266053812645353,0.2032838495769762,0.2130320280924107,-0.2957614816815313,0.0178354610241942,0.9274578856965077,0.8492910865982255,-0.2887541242905018,-0.9913772039381663,0.3216084305343259,-1.518989985892376,0.2656311150193504,0.5380972374482917,-0.7931664771795908,-1.774111517520497,0.8681779091439705,-0.8524941398257859,-0.6658269435512991,0.8796132488258497,-0.34072747683543,-1.987985903767425,0.5456950949477561,-1.782357498582343,0.05732155456410304,0.6009978402837941,0.07882144237710875,-0.8979680004924104,-1.430727923011459,0.5530271785698857,0.1107662425422351,-1.425308271745022,0.1700360757237642,-1.453625712786084,-1.169269120546788,-0.0126405038038242,-0.9594281116381027,-0.9739617500713154,-0.8530020464754762,-1.39827311402813,-0.9188941847089076,-1.386108329083921,0.7108487030822852,-0.8836487540671554,0.07853201812980926,-1.051325459232974,-1.484826849227233,-1.533234775666168,0.03140580681978467,0.02705202910173066,0.4166152682904976,0.3632893376491442,-1.420796830928111,-0.286231475378727,-0.3330877903531622,0.409372892150798,0.2198134769773095,-0.4647910620040745,-0.7310248159645814,-0.634855581722761,-0.6522933741681762,0.6972456050590643,-0.2740109892205576,-0.0131039015436536,-0.859563232963751,-1.706834688023223,-1.896388038219382,0.05557497110752241,-0.08979634413806226,0.2647219626799888,-1.419882546819699,-0.97020560221097,-0.2537171860121687,-1.147634038497515,-1.259516739674108,-1.995695558241755,0.1772621099422778,0.9216093755304966,-1.167901623284115,-0.4633035952839877,-1.103695748025646,-1.438126794358277,-0.7662231343711063,-1.382006963127418,0.5025553164634449,-1.708066538519708,-0.3926365821366558,0.5444420221363693,-1.497750219019922,0.9280863137415682,0.1756034011420033,0.65014857960709,-1.062815323873509,-0.6988142052810771,0.2901310367029266,0.4777884998235513,-1.035442613393293,0.8358822391478764)
for ($i = 0; $i -lt $uVals.Length; $i++) {
    $wsU.Cells.Item($i + 2, 2).Value = $uVals[$i]
}

# ---------------------------------------------------------------------------
# 4) "x" sheet (tab now named "x"): remove the old "i" header cell from A1
#    (content + formatting), relabel the value header in B1 as "x" and
#    write the new x-values.
# ---------------------------------------------------------------------------
$wsX.Range("A1").Clear()
$wsX.Range("B1").Value = "x"

$xVals = @(-0.008362533212884391,-0.008196812686071159,-0.008411903098170102,-0.00821109484127867,-0.008270119841687131,-0.008176634735496614,-0.008234458233781168,-0.008337811453716169,-0.008236806918498114,-0.008583333255869059,-0.008716426062069208,-0.008446915936975703,-0.008629675227311413,-0.008543805173801268,-0.008732897617185429,-0.008476897609003851,-0.008571471226388823,-0.00857343296194314,-0.008670833684618752,-0.008706457038653755,-0.008632814845064153,-0.008145599163935339,-0.008266143275096309,-0.008270385173598538,0.8175742996378883,-0.008669849181428581,-0.008387482276836415,0.9373030810195462,-0.008704037447086126,-0.008361850619191154,-0.008286355143636003,0.7988499232394671,-0.008395383061835274,-0.008148064894272739,-0.008353746977275145,-0.008749714703404785,-0.008379795398783103,-0.008263316394928617,-0.008698038659416946,-0.008430573289359557,-0.008636812963755791,-0.008310271005368644,-0.008174019075662906,-0.008715361584402607,-0.008285961109697119,-0.008278639074863205,-0.008849128051278521,-0.008725452550966392,-0.008610327860614482,-0.008369741005484869,-0.008507217725738796,-0.008697568900717786,-0.008560928908524272,0.7914252974303428,-0.008672860658859525,-0.008675511921058619,-0.008541326809775433,-0.008623036426752483,-0.008879195621726994,-0.008855970246774881,-0.008543119078057172,-0.008370527182817356,-0.008705263448752026,0.86714867199784,-0.008689873607786199,-0.008765834918330707,-0.008417806296953882,-0.008193491886609465,-0.008861568558451624,-0.008403550227774694,-0.008448757052192559,-0.008864972253482147,-0.008529864223059172,-0.008147623012643488,-0.9171013017100513,-0.00819170752025911,-0.008633550850652393,-0.9499783052364051,-0.008639296745550279,-0.008392685287222411,-0.0082693291362661,-0.008770082076552759,-0.008647860497294329,-0.008270554282885423,-0.008663843879103536,-0.008264192106702862,-0.008328967735009872,-0.008614960823736816,-0.008378075428502904,-0.008374641105153275,-0.008403428461618864,-0.008276645011612254,-0.008387698447860343,-0.008279389321743486,-0.9915140383047166,-0.00839610327676057,-0.00863921931735066,-0.008356428876496194,-0.008257204570087081,-0.008246406509202622,-0.008626644016284286,-0.00862548786519984,-0.008731614971605572,-0.008716795907616155,-0.008271567698101371,-0.008543766885705512,-0.008531808298432414,-0.008729614994057305,-0.008677358356134392,-0.00849854431954042,-0.008432856859328293,-0.008456353352311079,-0.008452073197607105,-0.9874424334866059,-0.008546897529972477,-0.008614838282860977,-0.008401840831310511,0.9035230168581602,0.9311335212955568,-0.008633084042683224,-0.008594640471990036,-0.008689625359364482,-0.008271768527650729,-0.0083755278026933,-0.008552105413073944,-0.008333978668599598,-0.0083081909603518,-0.008145975579373943,-0.008665799555364,-0.008877469037212361,-0.008329287573658689,-0.008498919787189262,-0.008344196396971431,-0.008267664307371957,-0.008424327939561035,-0.008280316331473943,-0.008755774667796101,-0.00820787179385893,-0.00851670610156869,-0.008767641201340546,-0.008254307093710615,-0.008879415177527267,-0.008665351273247449,-0.008797851951968945,-0.008353736613842973,-0.008440698406453442,-0.008696595653236973,-0.008748792207157647,-0.008360155105526122,-1.020577878880851)
for ($i = 0; $i -lt $xVals.Length; $i++) {
    $wsX.Cells.Item($i + 2, 2).Value = $xVals[$i]
}

# ---------------------------------------------------------------------------
# 5) "y" sheet: flip the indicated rows from 0 to 1.
# ---------------------------------------------------------------------------
$wsY = $wb.Worksheets.Item("y")
$yRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141,142,143,144,145,146,147,148,149,150,151)
foreach ($r in $yRows) {
    $wsY.Cells.Item($r, 2).Value = 1
}

Write-Host "edit complete"
